{"js": "// The badge table holds \"NameNN\" / \"InfoNN\" / \"ImageNN\" labels. For badges\n// 10-19 and 21, the \"Info\"/\"Image\" word is its own run, but the two-digit\n// number after it was originally typed as two separate single-character\n// runs (e.g. \"1\" then \"0\"). This merges just those two digit runs into a\n// single run containing the full two-digit number, leaving the\n// \"Info\"/\"Image\" run untouched - matching how the single-digit badges\n// (04-09, 20, 22-24) already store their number as one run.\nconst targets = [\n  { label: \"Info\", number: \"10\" },\n  { label: \"Image\", number: \"10\" },\n  { label: \"Info\", number: \"11\" },\n  { label: \"Image\", number: \"11\" },\n  { label: \"Info\", number: \"12\" },\n  { label: \"Image\", number: \"13\" }, // also matches the (pre-existing, duplicated) Image value for badge 13\n  { label: \"Info\", number: \"13\" },\n  { label: \"Info\", number: \"14\" },\n  { label: \"Image\", number: \"14\" },\n  { label: \"Info\", number: \"15\" },\n  { label: \"Image\", number: \"15\" },\n  { label: \"Info\", number: \"16\" },\n  { label: \"Image\", number: \"16\" },\n  { label: \"Info\", number: \"17\" },\n  { label: \"Image\", number: \"17\" },\n  { label: \"Info\", number: \"18\" },\n  { label: \"Image\", number: \"18\" },\n  { label: \"Info\", number: \"19\" },\n  { label: \"Image\", number: \"19\" },\n  { label: \"Info\", number: \"21\" },\n  { label: \"Image\", number: \"21\" }\n];\n\nfor (const { label, number } of targets) {\n  const fullText = label + number;\n  const matches = context.document.body.search(fullText, { matchCase: true, matchWildcards: false });\n  matches.load(\"items\");\n  await context.sync();\n\n  for (const matchRange of matches.items) {\n    // Narrow down to just the \"label\" part of this match, then take the\n    // range spanning from right after it to the end of the match - i.e.\n    // just the numeric suffix - and collapse that into one run.\n    const labelMatches = matchRange.search(label, { matchCase: true, matchWildcards: false });\n    labelMatches.load(\"items\");\n    await context.sync();\n\n    const labelRange = labelMatches.items[0];\n    const afterLabel = labelRange.getRange(\"After\");\n    const matchEnd = matchRange.getRange(\"End\");\n    const numberRange = afterLabel.expandTo(matchEnd);\n\n    numberRange.insertText(number, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# The badge table holds \"NameNN\" / \"InfoNN\" / \"ImageNN\" labels whose two-\n# digit numeric suffix (e.g. \"10\") was originally typed as two separate\n# runs (\"1\" then \"0\"). This merges the \"Info\"/\"Image\" runs for badges\n# 10-19 and 21 into a single run containing the full two-digit number,\n# matching how the single-digit badges (04-09, 20, 22-24) are already\n# stored as one run. Find/Replace (with identical find/replace text) is\n# used purely to force Word to collapse the run boundaries.\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Info10\", \"Image10\",\n  \"Info11\", \"Image11\",\n  \"Info12\",\n  \"Image13\", # also matches the (pre-existing, duplicated) Image value for badge 13\n  \"Info13\",\n  \"Info14\", \"Image14\",\n  \"Info15\", \"Image15\",\n  \"Info16\", \"Image16\",\n  \"Info17\", \"Image17\",\n  \"Info18\", \"Image18\",\n  \"Info19\", \"Image19\",\n  \"Info21\", \"Image21\"\n)\n\nforeach ($target in $targets) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)\n}\n"}
